$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# The homework numbering shifts by one starting with item 2: a brand new
# question ("2. Import the file sleep.txt into SPSS. ...") is inserted in
# place of the old item 2, the old item-2 text ("Fit a multiple linear
# regression model ...") becomes the new item 3, and every item from the old
# 3 through 7 is bumped up by one (3->4, 4->5, 5->6, 6->7, 7->8).
#
# We first bump the numeral prefixes of the later, untouched-in-content
# paragraphs (working from the bottom up, using their original paragraph
# indices) and only then rewrite/insert the paragraphs around the old item 2,
# since that insertion is what shifts every later paragraph index.
# ---------------------------------------------------------------------------

function Bump-LeadingNumber($paraIndex, $newDigit) {
    $p = $d.Paragraphs.Item($paraIndex)
    $rng = $d.Range($p.Range.Start, $p.Range.Start + 1)
    $rng.Text = $newDigit
}

# item 7 -> 8  ("7. Compute a correlation matrix ...")
Bump-LeadingNumber 16 "8"

# item 6 -> 7  ("6. Draw a scatterplot ...")
Bump-LeadingNumber 14 "7"

# item 5 -> 6  ("5. Calculate the residuals ...")
Bump-LeadingNumber 12 "6"

# item 4 -> 5  ("4. Show the value of R-squared ...")
Bump-LeadingNumber 10 "5"

# item 3 -> 4  ("3. Ignore the reviewers objection ...")
Bump-LeadingNumber 8 "4"

# ---------------------------------------------------------------------------
# Now handle the old item-2 paragraph: its question text is replaced by the
# new SPSS/sleep.txt import question, a blank paragraph is inserted after it,
# and the original question text re-appears (renumbered "3.") in a brand new
# paragraph right after that blank one.
# ---------------------------------------------------------------------------

$p2 = $d.Paragraphs.Item(6)
$fullText = $p2.Range.Text
$fullText = $fullText.Substring(0, $fullText.Length - 1)   # drop trailing paragraph mark
$oldQuestionBody = $fullText.Substring(3)                  # drop leading "2. "

$p2Body = $d.Range($p2.Range.Start, $p2.Range.End - 1)
$p2Body.Text = "2. Import the file sleep.txt into SPSS. Refer to the data dictionary if needed. Some of the variables may be misidentified as strings rather than numeric, so please check your data carefully after importing it. Display the first ten rows of data below."

$p2 = $d.Paragraphs.Item(6)
$p2.Range.InsertParagraphAfter()

$pBlank = $d.Paragraphs.Item(7)
$pBlank.Range.InsertParagraphAfter()

$pItem3 = $d.Paragraphs.Item(8)
$pItem3Body = $d.Range($pItem3.Range.Start, $pItem3.Range.End - 1)
$pItem3Body.Text = "3. " + $oldQuestionBody
